$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '67.614.47'
$ws.Range('D3').Value = '3.783.69'
$ws.Range('E3').Value = '  +0.94%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''594.91'
$ws.Range('E5').Value = '  +0.26%  '
$ws.Range('D6').Value = '''166.89'
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').Value = '3.769.42'
$ws.Range('E7').Value = '  +0.57%  '
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('E10').Value = '  -0.38%  '
$ws.Range('D12').Value = '''0.448'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('E13').Value = '  -2.78%  '
$ws.Range('D14').Value = '''35.97'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').Value = '4.419.06'
$ws.Range('E15').Value = '  +1.00%  '
$ws.Range('D16').Value = '3.779.40'
$ws.Range('E16').Value = '  +0.77%  '
$ws.Range('E17').Value = '  +3.62%  '
$ws.Range('D18').Value = '67.602.85'
$ws.Range('E18').Value = '  -1.19%  '
$ws.Range('D19').Value = '''7.01'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('E20').Value = '  -0.05%  '
$ws.Range('D21').Value = '''10.03'
$ws.Range('E21').Value = '  -5.77%  '
$ws.Range('D22').Value = '''459.11'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('D23').Value = '''0.695'
$ws.Range('E23').Value = '  -0.32%  '
$ws.Range('D24').Value = '''0.0000154'
$ws.Range('E24').Value = '  +5.19%  '
$ws.Range('D25').Value = '''83.36'
$ws.Range('E25').Value = '  -0.78%  '
$ws.Range('D26').Value = '''11.98'
$ws.Range('E26').Value = '  +0.77%  '
$ws.Range('E27').Value = '  -3.03%  '
$ws.Range('E28').Value = '  +0.12%  '
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = '''2.22'
$ws.Range('E31').Value = '  +3.04%  '
$ws.Range('E32').Value = '  -1.33%  '
$ws.Range('D33').Value = '''29.58'
$ws.Range('E33').Value = '  -0.91%  '
$ws.Range('D34').Value = '''0.998'
$ws.Range('E34').Value = '  -0.26%  '
$ws.Range('D35').Value = '''9.07'
$ws.Range('E35').Value = '  -1.11%  '
$ws.Range('E36').Value = '  -0.51%  '
$ws.Range('D37').Value = '''3.36'
$ws.Range('E37').Value = '  -0.73%  '
$ws.Range('E38').Value = '  -0.23%  '
$ws.Range('D39').Value = '''0.995'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').Value = '''5.75'
$ws.Range('E40').Value = '  -0.50%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D43').Value = '''45.22'
$ws.Range('E43').Value = '  +3.28%  '
$ws.Range('D44').Value = '''48.14'
$ws.Range('E44').Value = '  +3.41%  '
$ws.Range('D45').Value = '''0.298'
$ws.Range('E45').Value = '  -0.93%  '
$ws.Range('E46').Value = '  +3.74%  '
$ws.Range('D47').Value = '''8.30'
$ws.Range('E47').Value = '  -1.85%  '
$ws.Range('D48').Value = '''394.65'
$ws.Range('E48').Value = '  +0.82%  '
$ws.Range('D49').Value = '''26.59'
$ws.Range('E49').Value = '  +5.88%  '
$ws.Range('E50').Value = '  -5.16%  '
$ws.Range('D51').Value = '2.715.34'
$ws.Range('E51').Value = '  -1.29%  '
